# Generate Report for Handback
#
# 1. Status text moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    on every sheet (Overview + each locale sheet).
# 2. Each locale sheet (zh-cn, de-de) gains a populated "Latest Target File" /
#    "Latest Handback File" / "Latest Handback DateTime" for the two real rows
#    (the .localization-config row stays untouched/"Ignored").

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$mdFile  = "1abafd1a-6f17-4f05-b1f0-8738f5748055.md"
$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/3509e5450a3f995873eefd3915d5b72949aceb27/e2e/1abafd1a-6f17-4f05-b1f0-8738f5748055.md"

$zhXlfFile = "1abafd1a-6f17-4f05-b1f0-8738f5748055.70d6c5d2cd32aba90dc8c64c0284a597814f2ea9.zh-cn.xlf"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d0517585a78df41890c24773ba0050a1a71efd68/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1abafd1a-6f17-4f05-b1f0-8738f5748055.70d6c5d2cd32aba90dc8c64c0284a597814f2ea9.zh-cn.xlf"
$zhHandbackDate = "2016-03-09 16:55:26"

$deXlfFile = "1abafd1a-6f17-4f05-b1f0-8738f5748055.70d6c5d2cd32aba90dc8c64c0284a597814f2ea9.de-de.xlf"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/efb86b8ddf045403b61983213d7c411aba8f0634/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1abafd1a-6f17-4f05-b1f0-8738f5748055.70d6c5d2cd32aba90dc8c64c0284a597814f2ea9.de-de.xlf"
$deHandbackDate = "2016-03-09 16:55:34"

# ---- Overview sheet: refresh the rolled-up status column ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B2").Value = $newStatus
$zh.Range("B3").Value = $newStatus

$zh.Range("E2").Value = $mdFile
$zh.Hyperlinks.Add($zh.Range("E2"), $mdUrl, "", "", $mdFile)
$zh.Range("F2").Value = $zhXlfFile
$zh.Hyperlinks.Add($zh.Range("F2"), $zhXlfUrl, "", "", $zhXlfFile)
$zh.Range("G2").Value = $zhHandbackDate

$zh.Range("E3").Value = $mdFile
$zh.Hyperlinks.Add($zh.Range("E3"), $mdUrl, "", "", $mdFile)
$zh.Range("F3").Value = $zhXlfFile
$zh.Hyperlinks.Add($zh.Range("F3"), $zhXlfUrl, "", "", $zhXlfFile)
$zh.Range("G3").Value = $zhHandbackDate

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("B2").Value = $newStatus
$de.Range("B3").Value = $newStatus

$de.Range("E2").Value = $mdFile
$de.Hyperlinks.Add($de.Range("E2"), $mdUrl, "", "", $mdFile)
$de.Range("F2").Value = $deXlfFile
$de.Hyperlinks.Add($de.Range("F2"), $deXlfUrl, "", "", $deXlfFile)
$de.Range("G2").Value = $deHandbackDate

$de.Range("E3").Value = $mdFile
$de.Hyperlinks.Add($de.Range("E3"), $mdUrl, "", "", $mdFile)
$de.Range("F3").Value = $deXlfFile
$de.Hyperlinks.Add($de.Range("F3"), $deXlfUrl, "", "", $deXlfFile)
$de.Range("G3").Value = $deHandbackDate
